# Refresh the crypto price/volume table with the latest scraped values.
# (Two rows also swapped rank order: WrappedBTC/TRON at 16-17 and
#  Mantle/ImmutableX at 31-32.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '''62.875.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.84%  '
$ws.Range("E2").Style = "Normal"
# Row 3: Ethereum
$ws.Range("D3").Value = '''3.067.39'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.99%  '
$ws.Range("E3").Style = "Normal"
# Row 4: TetherUSD
$ws.Range("E4").Value = '''  -0.34%  '
$ws.Range("E4").Style = "Normal"
# Row 5: BNB
$ws.Range("D5").Value = '''541.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.42%  '
$ws.Range("E5").Style = "Normal"
# Row 6: Solana
$ws.Range("D6").Value = '''137.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +3.93%  '
$ws.Range("E6").Style = "Normal"
# Row 7: USDC
$ws.Range("E7").Value = '''  -0.25%  '
$ws.Range("E7").Style = "Normal"
# Row 8: LidoStakedEther
$ws.Range("D8").Value = '''3.061.97'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.94%  '
$ws.Range("E8").Style = "Normal"
# Row 9: XRP
$ws.Range("E9").Value = '''  +1.91%  '
$ws.Range("E9").Style = "Normal"
# Row 10: Dogecoin
$ws.Range("E10").Value = '''  +2.55%  '
$ws.Range("E10").Style = "Normal"
# Row 11: Toncoin
$ws.Range("D11").Value = '''6.27'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +2.81%  '
$ws.Range("E11").Style = "Normal"
# Row 12: Cardano
$ws.Range("D12").Value = '''0.456'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -0.48%  '
$ws.Range("E12").Style = "Normal"
# Row 13: ShibaInu
$ws.Range("E13").Value = '''  +5.69%  '
$ws.Range("E13").Style = "Normal"
# Row 14: Avalanche
$ws.Range("D14").Value = '''34.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +0.72%  '
$ws.Range("E14").Style = "Normal"
# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '''3.566.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +0.73%  '
$ws.Range("E15").Style = "Normal"
# Row 16: TRON
$ws.Range("B16").Value = '''TRON'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = '''https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '''0.113'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +1.89%  '
$ws.Range("E16").Style = "Normal"
# Row 17: WrappedBTC
$ws.Range("B17").Value = '''WrappedBTC'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '''62.909.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +0.45%  '
$ws.Range("E17").Style = "Normal"
# Row 18: WrappedEther
$ws.Range("D18").Value = '''3.070.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.11%  '
$ws.Range("E18").Style = "Normal"
# Row 19: Polkadot
$ws.Range("D19").Value = '''6.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +1.05%  '
$ws.Range("E19").Style = "Normal"
# Row 20: BitcoinCash
$ws.Range("D20").Value = '''471.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -0.92%  '
$ws.Range("E20").Style = "Normal"
# Row 21: Chainlink
$ws.Range("D21").Value = '''13.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +2.15%  '
$ws.Range("E21").Style = "Normal"
# Row 22: Polygon
$ws.Range("E22").Value = '''  -0.22%  '
$ws.Range("E22").Style = "Normal"
# Row 23: Uniswap
$ws.Range("D23").Value = '''7.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -0.41%  '
$ws.Range("E23").Style = "Normal"
# Row 24: Litecoin
$ws.Range("D24").Value = '''78.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.34%  '
$ws.Range("E24").Style = "Normal"
# Row 25: InternetComputer(DFINITY)
$ws.Range("E25").Value = '''  +0.87%  '
$ws.Range("E25").Style = "Normal"
# Row 26: Dai
$ws.Range("E26").Value = '''  +0.48%  '
$ws.Range("E26").Style = "Normal"
# Row 27: PancakeSwap
$ws.Range("E27").Value = '''  +0.13%  '
$ws.Range("E27").Style = "Normal"
# Row 28: RenderToken
$ws.Range("D28").Value = '''7.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -4.19%  '
$ws.Range("E28").Style = "Normal"
# Row 29: FirstDigitalUSD
$ws.Range("E29").Value = '''  -0.09%  '
$ws.Range("E29").Style = "Normal"
# Row 30: EthereumClassic
$ws.Range("D30").Value = '''26.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +0.49%  '
$ws.Range("E30").Style = "Normal"
# Row 31: ImmutableX
$ws.Range("B31").Value = '''ImmutableX'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = '''1.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -0.67%  '
$ws.Range("E31").Style = "Normal"
# Row 32: Mantle
$ws.Range("B32").Value = '''Mantle'
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = '''1.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +5.27%  '
$ws.Range("E32").Style = "Normal"
# Row 33: OKB
$ws.Range("E33").Value = '''  -2.29%  '
$ws.Range("E33").Style = "Normal"
# Row 34: NEARProtocol
$ws.Range("D34").Value = '''5.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +11.06%  '
$ws.Range("E34").Style = "Normal"
# Row 35: Stacks
$ws.Range("D35").Value = '''2.30'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -5.26%  '
$ws.Range("E35").Style = "Normal"
# Row 36: Filecoin
$ws.Range("D36").Value = '''5.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +1.77%  '
$ws.Range("E36").Style = "Normal"
# Row 37: Bittensor
$ws.Range("D37").Value = '''485.57'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -1.20%  '
$ws.Range("E37").Style = "Normal"
# Row 38: Maker
$ws.Range("D38").Value = '''3.256.63'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +4.58%  '
$ws.Range("E38").Style = "Normal"
# Row 39: VeChain
$ws.Range("E39").Value = '''  +2.65%  '
$ws.Range("E39").Style = "Normal"
# Row 40: Hedera
$ws.Range("D40").Value = '''0.0793'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +1.38%  '
$ws.Range("E40").Style = "Normal"
# Row 41: Kaspa
$ws.Range("D41").Value = '''0.118'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +1.17%  '
$ws.Range("E41").Style = "Normal"
# Row 42: Cosmos
$ws.Range("D42").Value = '''8.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +2.43%  '
$ws.Range("E42").Style = "Normal"
# Row 43: dogwifhat
$ws.Range("D43").Value = '''2.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +2.20%  '
$ws.Range("E43").Style = "Normal"
# Row 44: TheGraph
$ws.Range("D44").Value = '''0.252'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +1.42%  '
$ws.Range("E44").Style = "Normal"
# Row 46: InjectiveProtocol
$ws.Range("D46").Value = '''25.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +2.85%  '
$ws.Range("E46").Style = "Normal"
# Row 47: Monero
$ws.Range("D47").Value = '''122.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +4.42%  '
$ws.Range("E47").Style = "Normal"
# Row 48: Fetch.AI
$ws.Range("D48").Value = '''2.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +0.46%  '
$ws.Range("E48").Style = "Normal"
# Row 49: Stellar
$ws.Range("E49").Value = '''  +2.53%  '
$ws.Range("E49").Style = "Normal"
# Row 50: PEPE
$ws.Range("D50").Value = '''0.0₃0523'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +5.45%  '
$ws.Range("E50").Style = "Normal"
# Row 51: ThetaToken
$ws.Range("D51").Value = '''2.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +1.81%  '
$ws.Range("E51").Style = "Normal"
